# Adjustment CUM/KWh - PJC
#
# Adds a new "Usage Type" label to the Raw Data template sheet, placed in
# cell A9 right below the existing "Building Description" label (A8), and
# styled to match the other bold header/label cells already on the sheet.
# Finally, the active selection is moved onto B9 (the value cell next to
# the new label), matching the saved worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A9")
$cell.Value = "Usage Type"
$cell.Font.Bold = $true

$null = $ws.Range("B9").Select()
